$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Fixed the git repo" bullet: the word "git" was wrapped in
#    proofErr spell-check markers, splitting the sentence across three
#    runs. Re-typing the sentence via Find/Replace collapses it back
#    into a single run and drops the now-stale proofErr elements.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Fixed the git repo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Fixed the git repo", 2)

# ---------------------------------------------------------------------
# 2. Remove the stray empty bullet (numId 7) that was left sitting
#    right after "Fix reference signal".
# ---------------------------------------------------------------------
$found = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Fix reference signal`r") {
        $found = $true
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -eq "`r") {
            $next.Range.Delete()
        }
        break
    }
}

Write-Output "done found=$found"
